$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 2
    3 = 1
    4 = 3
    5 = 1
    6 = 3
    7 = 2
    8 = 1
    9 = 2
    10 = 4
    11 = 3
    12 = 1
    13 = 6
    14 = 1
    15 = 2
    16 = 4
    17 = 1
    18 = 3
    19 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
